$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# NOTE: this engine's Find.Execute matches across the whole document body
# rather than confining the search to the receiver Range, even though the
# replacement is applied only to the single match that Find locates (it
# does not respect wdReplaceAll scoping to the range either). To keep each
# replacement unambiguous we (a) target the specific table cell so unique
# text still resolves correctly, and (b) order the 25 edits so that no two
# cells ever hold identical text at the moment a Find runs - in particular
# Cell(13,1) "570÷6=" -> "828÷5=" must happen BEFORE Cell(5,2)
# "360÷8=" -> "570÷6=" creates a second, temporary "570÷6=".
$found = $t.Cell(1,1).Range.Find.Execute("898÷7=", $true, $false, $false, $false, $false, $true, 0, $false, "891÷8=", 1)
if (-not $found) { Write-Host "WARNING: replace failed for Cell(1,1): 898÷7= -> 891÷8=" }
$found = $t.Cell(1,2).Range.Find.Execute("674÷5=", $true, $false, $false, $false, $false, $true, 0, $false, "378÷2=", 1)
if (-not $found) { Write-Host "WARNING: replace failed for Cell(1,2): 674÷5= -> 378÷2=" }
$found = $t.Cell(1,3).Range.Find.Execute("377÷5=", $true, $false, $false, $false, $false, $true, 0, $false, "946÷8=", 1)
if (-not $found) { Write-Host "WARNING: replace failed for Cell(1,3): 377÷5= -> 946÷8=" }
$found = $t.Cell(1,4).Range.Find.Execute("412÷2=", $true, $false, $false, $false, $false, $true, 0, $false, "368÷6=", 1)
if (-not $found) { Write-Host "WARNING: replace failed for Cell(1,4): 412÷2= -> 368÷6=" }
$found = $t.Cell(1,5).Range.Find.Execute("920÷8=", $true, $false, $false, $false, $false, $true, 0, $false, "570÷5=", 1)
if (-not $found) { Write-Host "WARNING: replace failed for Cell(1,5): 920÷8= -> 570÷5=" }
$found = $t.Cell(5,1).Range.Find.Execute("332÷9=", $true, $false, $false, $false, $false, $true, 0, $false, "452÷9=", 1)
if (-not $found) { Write-Host "WARNING: replace failed for Cell(5,1): 332÷9= -> 452÷9=" }
$found = $t.Cell(13,1).Range.Find.Execute("570÷6=", $true, $false, $false, $false, $false, $true, 0, $false, "828÷5=", 1)
if (-not $found) { Write-Host "WARNING: replace failed for Cell(13,1): 570÷6= -> 828÷5=" }
$found = $t.Cell(5,2).Range.Find.Execute("360÷8=", $true, $false, $false, $false, $false, $true, 0, $false, "570÷6=", 1)
if (-not $found) { Write-Host "WARNING: replace failed for Cell(5,2): 360÷8= -> 570÷6=" }
$found = $t.Cell(5,3).Range.Find.Execute("277÷8=", $true, $false, $false, $false, $false, $true, 0, $false, "146÷5=", 1)
if (-not $found) { Write-Host "WARNING: replace failed for Cell(5,3): 277÷8= -> 146÷5=" }
$found = $t.Cell(5,4).Range.Find.Execute("376÷2=", $true, $false, $false, $false, $false, $true, 0, $false, "470÷6=", 1)
if (-not $found) { Write-Host "WARNING: replace failed for Cell(5,4): 376÷2= -> 470÷6=" }
$found = $t.Cell(5,5).Range.Find.Execute("198÷9=", $true, $false, $false, $false, $false, $true, 0, $false, "405÷3=", 1)
if (-not $found) { Write-Host "WARNING: replace failed for Cell(5,5): 198÷9= -> 405÷3=" }
$found = $t.Cell(9,1).Range.Find.Execute("887÷5=", $true, $false, $false, $false, $false, $true, 0, $false, "684÷3=", 1)
if (-not $found) { Write-Host "WARNING: replace failed for Cell(9,1): 887÷5= -> 684÷3=" }
$found = $t.Cell(9,2).Range.Find.Execute("160÷6=", $true, $false, $false, $false, $false, $true, 0, $false, "975÷4=", 1)
if (-not $found) { Write-Host "WARNING: replace failed for Cell(9,2): 160÷6= -> 975÷4=" }
$found = $t.Cell(9,3).Range.Find.Execute("562÷3=", $true, $false, $false, $false, $false, $true, 0, $false, "954÷6=", 1)
if (-not $found) { Write-Host "WARNING: replace failed for Cell(9,3): 562÷3= -> 954÷6=" }
$found = $t.Cell(9,4).Range.Find.Execute("326÷6=", $true, $false, $false, $false, $false, $true, 0, $false, "300÷9=", 1)
if (-not $found) { Write-Host "WARNING: replace failed for Cell(9,4): 326÷6= -> 300÷9=" }
$found = $t.Cell(9,5).Range.Find.Execute("127÷7=", $true, $false, $false, $false, $false, $true, 0, $false, "671÷5=", 1)
if (-not $found) { Write-Host "WARNING: replace failed for Cell(9,5): 127÷7= -> 671÷5=" }
$found = $t.Cell(13,2).Range.Find.Execute("651÷4=", $true, $false, $false, $false, $false, $true, 0, $false, "798÷3=", 1)
if (-not $found) { Write-Host "WARNING: replace failed for Cell(13,2): 651÷4= -> 798÷3=" }
$found = $t.Cell(13,3).Range.Find.Execute("100÷2=", $true, $false, $false, $false, $false, $true, 0, $false, "824÷8=", 1)
if (-not $found) { Write-Host "WARNING: replace failed for Cell(13,3): 100÷2= -> 824÷8=" }
$found = $t.Cell(13,4).Range.Find.Execute("289÷2=", $true, $false, $false, $false, $false, $true, 0, $false, "717÷8=", 1)
if (-not $found) { Write-Host "WARNING: replace failed for Cell(13,4): 289÷2= -> 717÷8=" }
$found = $t.Cell(13,5).Range.Find.Execute("649÷6=", $true, $false, $false, $false, $false, $true, 0, $false, "681÷6=", 1)
if (-not $found) { Write-Host "WARNING: replace failed for Cell(13,5): 649÷6= -> 681÷6=" }
$found = $t.Cell(17,1).Range.Find.Execute("817÷9=", $true, $false, $false, $false, $false, $true, 0, $false, "866÷2=", 1)
if (-not $found) { Write-Host "WARNING: replace failed for Cell(17,1): 817÷9= -> 866÷2=" }
$found = $t.Cell(17,2).Range.Find.Execute("145÷2=", $true, $false, $false, $false, $false, $true, 0, $false, "496÷9=", 1)
if (-not $found) { Write-Host "WARNING: replace failed for Cell(17,2): 145÷2= -> 496÷9=" }
$found = $t.Cell(17,3).Range.Find.Execute("652÷4=", $true, $false, $false, $false, $false, $true, 0, $false, "988÷6=", 1)
if (-not $found) { Write-Host "WARNING: replace failed for Cell(17,3): 652÷4= -> 988÷6=" }
$found = $t.Cell(17,4).Range.Find.Execute("363÷7=", $true, $false, $false, $false, $false, $true, 0, $false, "178÷6=", 1)
if (-not $found) { Write-Host "WARNING: replace failed for Cell(17,4): 363÷7= -> 178÷6=" }
$found = $t.Cell(17,5).Range.Find.Execute("652÷2=", $true, $false, $false, $false, $false, $true, 0, $false, "119÷4=", 1)
if (-not $found) { Write-Host "WARNING: replace failed for Cell(17,5): 652÷2= -> 119÷4=" }

Write-Host "Done."
